# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 59
$ws1.Range("F8").Value = 86
$ws1.Range("F9").Value = 8557
$ws1.Range("F10").Value = 797
$ws1.Range("F11").Value = 319
$ws1.Range("F12").Value = 1137
$ws1.Range("F13").Value = 931
$ws1.Range("F14").Value = 92
$ws1.Range("F17").Value = 214
$ws1.Range("F19").Value = 228
$ws1.Range("F20").Value = 975

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 59
$ws4.Range("F10").Value = 86
$ws4.Range("F11").Value = 8557
$ws4.Range("F12").Value = 797
$ws4.Range("F13").Value = 319
$ws4.Range("F14").Value = 1137
$ws4.Range("F15").Value = 931
$ws4.Range("F16").Value = 92
$ws4.Range("F19").Value = 214
$ws4.Range("F21").Value = 228
$ws4.Range("F22").Value = 975
